$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Insert one new row at row 44 (shifts existing rows 44..187 down to
# 45..188, and extends the used range to A1:J188).
# ---------------------------------------------------------------------
$ws.Rows.Item(44).Insert() | Out-Null

# The emulator's row-insert does not faithfully clone the neighboring
# row's cell formatting, so explicitly reuse the formats already present
# on row 43 (the row now directly above the new blank row) via a
# formats-only paste. This keeps the existing style records (border-only
# style + bordered/2-decimal numeric style) instead of minting new ones.
$ws.Range("A43:D43").Copy() | Out-Null
$ws.Range("A44:D44").PasteSpecial(-4122) | Out-Null

$ws.Range("E43:F43").Copy() | Out-Null
$ws.Range("E44:F44").PasteSpecial(-4122) | Out-Null

$ws.Range("H43").Copy() | Out-Null
$ws.Range("H44").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Populate the newly inserted row 44 with the new product entry.
# ---------------------------------------------------------------------
$ws.Range("A44").Value = 6010000223
$ws.Range("B44").Value = "Clavmycin 625 DuoTab  (Alu-Alu) (New)"
$ws.Range("C44").Value = "20X1X10"
$ws.Range("D44").Value = 360
$ws.Range("E44").Value = 204.96
$ws.Range("F44").Value = 45.9
$ws.Range("H44").Value = 45.9

# ---------------------------------------------------------------------
# Match the author's resulting view state: scrolled so row 25 is the
# top-left visible row, with B44 as the active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("B44").Select() | Out-Null
